$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Date" column (B) holds a mix of real dates and text-formatted dates
# ("2026-03-21" typed in as text). Normalize B2:B4 into real Excel date
# values and apply the same date number format already used by B5
# (m/d/yyyy) to the whole date column, including the blank B6 cell.

$ws.Range("B2").Value = 46102
$ws.Range("B3").Value = 46102
$ws.Range("B4").Value = 46102

$ws.Range("B5").Copy()
$ws.Range("B2:B4").PasteSpecial(-4122)
$ws.Range("B6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Let Excel recompute the (now-uniform) row heights instead of the stale
# explicit heights left over from the previous formatting pass.
$ws.Rows("1:6").AutoFit()

# Tidy the view to match the saved state.
$ws.Range("B6").Select() | Out-Null
